$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "UseCase" -> "Use Case" (splits the quoted term into "Use" / " " / "Case"
#    and removes the spell-check proofErr markers that wrapped "UseCase").
# ---------------------------------------------------------------------------
$oldUseCase = [char]0x201C + "UseCase" + [char]0x201D
$newUseCase = [char]0x201C + "Use Case" + [char]0x201D
$r1 = $d.Content
$found1 = $r1.Find.Execute($oldUseCase, $false, $false, $false, $false, $false, $true, 1, $false, $newUseCase, 2)

# ---------------------------------------------------------------------------
# 2) "pianificazione, in seguito" -> "pianificazione. In seguito"
# ---------------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute(", in seguito ho usato", $false, $false, $false, $false, $false, $true, 1, $false, ". In seguito ho usato", 2)

# ---------------------------------------------------------------------------
# 3) Re-create the run boundaries that Word shows in the saved package by
#    toggling (and immediately un-toggling) a character attribute at each
#    split point - this forces the run to break at that offset without
#    altering the visible formatting.
# ---------------------------------------------------------------------------
$docEnd = $d.Content.End

function SplitAt([int]$pos) {
    $s = $d.Range($pos, $docEnd)
    $s.Bold = 1
    $s.Bold = 0
}

# Locate "\u201cUse" / " " / "Case" boundaries.
$rQ = $d.Content
$rQ.Find.Execute([char]0x201C) | Out-Null
$quoteStart = $rQ.Start

$rUC = $d.Content
$rUC.Find.Execute("Use Case") | Out-Null
$useCaseStart = $rUC.Start

SplitAt $quoteStart          # before "\u201cUse"
SplitAt ($useCaseStart + 3)  # before " " (after "Use")
SplitAt ($useCaseStart + 4)  # before "Case"

# Locate the "pianificazione. I" / "n seguito" boundaries.
$rP = $d.Content
$rP.Find.Execute("pianificazione. In seguito") | Out-Null
$pianifStart = $rP.Start   # start of "pianificazione"

$periodPos = $pianifStart + 14   # "pianificazione" is 14 chars long
SplitAt $periodPos            # before "."
SplitAt ($periodPos + 1)      # before " "
SplitAt ($periodPos + 2)      # before "I"
$bookmarkPos = $periodPos + 3 # before "n seguito..."

# ---------------------------------------------------------------------------
# 4) Move the "_GoBack" bookmark so it sits right before "n seguito..."
#    (Word keeps a single "_GoBack" bookmark, so adding the new one also
#    removes the stale one that used to sit at the end of the document).
# ---------------------------------------------------------------------------
$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
